$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from Q1 into the new R1 header cell
$ws.Range("Q1").Copy($ws.Range("R1"))
$ws.Range("R1").Value = "genre"

$ws.Range("R2").Value = "['canadian pop', 'pop', 'post-teen pop', 'viral pop']"
$ws.Range("R3").Value = "['latin', 'reggaeton flow', 'trap latino']"
$ws.Range("R4").Value = "['dance pop', 'pop', 'post-teen pop']"
$ws.Range("R5").Value = "['pop', 'uk pop']"
$ws.Range("R6").Value = "['dfw rap', 'melodic rap', 'rap']"
$ws.Range("R7").Value = "['pop', 'uk pop']"
$ws.Range("R8").Value = "['melodic rap']"
$ws.Range("R9").Value = "['pop', 'uk pop']"
$ws.Range("R10").Value = "['country rap', 'lgbtq+ hip hop', 'pop rap']"
$ws.Range("R11").Value = "['electropop', 'pop']"
$ws.Range("R12").Value = "['latin', 'reggaeton', 'trap latino']"
$ws.Range("R13").Value = "['electronic trap']"
$ws.Range("R14").Value = "['pop', 'uk pop']"
$ws.Range("R15").Value = "['panamanian pop', 'reggaeton']"
$ws.Range("R16").Value = "['canadian hip hop', 'canadian pop', 'hip hop', 'pop rap', 'rap', 'toronto rap']"
$ws.Range("R17").Value = "['dance pop', 'pop', 'pop rap', 'r&b', 'rap']"
$ws.Range("R18").Value = "['latin', 'reggaeton']"
$ws.Range("R19").Value = "['dfw rap', 'melodic rap', 'rap']"
$ws.Range("R20").Value = "['pop rap']"
$ws.Range("R21").Value = "['escape room', 'minnesota hip hop', 'pop', 'pop rap', 'trap queen']"
$ws.Range("R22").Value = "['pop house']"
$ws.Range("R23").Value = "['country rap', 'lgbtq+ hip hop', 'pop rap']"
$ws.Range("R24").Value = "['latin', 'reggaeton', 'reggaeton flow', 'trap latino']"
$ws.Range("R25").Value = "['latin', 'reggaeton', 'reggaeton flow']"
$ws.Range("R26").Value = "['electropop', 'pop']"
$ws.Range("R27").Value = "['canadian pop', 'pop', 'post-teen pop', 'viral pop']"
$ws.Range("R28").Value = "['australian pop']"
$ws.Range("R29").Value = "['canadian hip hop', 'pop']"
$ws.Range("R30").Value = "['latin', 'latin hip hop', 'reggaeton', 'tropical']"
$ws.Range("R31").Value = "['latin', 'reggaeton']"
$ws.Range("R32").Value = "['electropop', 'pop', 'tropical house']"
$ws.Range("R33").Value = "['dance pop', 'pop', 'post-teen pop']"
$ws.Range("R34").Value = "['latin', 'reggaeton']"
$ws.Range("R35").Value = "['atl hip hop', 'atl trap', 'gangster rap', 'melodic rap', 'pop rap', 'rap', 'trap']"
$ws.Range("R36").Value = "['dance pop', 'pop', 'post-teen pop']"
$ws.Range("R37").Value = "['big room', 'edm', 'pop', 'progressive house', 'tropical house']"
$ws.Range("R38").Value = "['panamanian pop', 'reggaeton']"
$ws.Range("R39").Value = "['pop', 'uk pop']"
$ws.Range("R40").Value = "['boy band', 'dance pop', 'pop', 'pop rock', 'post-teen pop']"
$ws.Range("R41").Value = "['pop']"
$ws.Range("R42").Value = "['edm', 'pop', 'tropical house']"
$ws.Range("R43").Value = "['pop', 'post-teen pop']"
$ws.Range("R44").Value = "['dance pop', 'pop']"
$ws.Range("R45").Value = "['alternative r&b', 'pop']"
$ws.Range("R46").Value = "['r&b en espanol']"
$ws.Range("R47").Value = "['brostep', 'progressive electro house']"
$ws.Range("R48").Value = "['latin', 'latin hip hop', 'reggaeton', 'tropical']"
$ws.Range("R49").Value = "['brostep', 'progressive electro house']"
$ws.Range("R50").Value = "['electropop', 'pop', 'tropical house']"
$ws.Range("R51").Value = "['pop', 'uk pop']"
